$wb = $excel.ActiveWorkbook

# --- Sheet "Dados": add the new "Lançado" status column (H) ---
$dados = $wb.Worksheets.Item("Dados")

# Copy the header style from the last existing header cell (G1) onto the new H1 header
$dados.Range("G1").Copy()
$dados.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$dados.Application.CutCopyMode = $false
$dados.Range("H1").Value = "Lançado"

# Add a list data-validation ("Pendente" / "Lançado") to the whole column
$statusRange = $dados.Range("H1:H1048576")
$statusRange.Validation.Add(3, 1, 1, """Pendente, Lançado""")
$statusRange.Validation.IgnoreBlank = $true
$statusRange.Validation.InCellDropdown = $true
$statusRange.Validation.ShowInput = $true
$statusRange.Validation.ShowError = $true

# --- Sheet "Atividade": drop the (unused) interior borders on columns B/C ---
$atividade = $wb.Worksheets.Item("Atividade")
$atividade.Range("B2:C10").Borders.LineStyle = -4142  # xlLineStyleNone
$atividade.Range("C2:C10").ClearContents()
$atividade.Range("B2").ClearContents()

Write-Output "edit applied"
